$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 11
$scratchStart = 100

# Step 1: copy the current data rows into a scratch area (well below the
# used range) so cell *type* (text vs number) and exact stored text is
# preserved verbatim - a plain Copy, not a value read/re-entry.
$n = $lastDataRow - $firstDataRow + 1
for ($i = 0; $i -lt $n; $i++) {
    $srcRow = $firstDataRow + $i
    $dstRow = $scratchStart + $i
    $src = $ws.Range("A" + $srcRow + ":D" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":D" + $dstRow)
    $src.Copy($dst)
}

# Step 2: determine the sort order of the scratch rows by their POINT name
# (column A), ascending, ordinary text sort.
$keys = @()
for ($i = 0; $i -lt $n; $i++) {
    $r = $scratchStart + $i
    $name = $ws.Cells.Item($r, 1).Value2
    $keys += ,@($name, $r)
}
$sortedKeys = $keys | Sort-Object { $_[0] }

# Step 3: move (Cut, not Copy) each scratch row into its sorted destination
# row in the original table - Cut preserves the cell's original type/text
# exactly, unlike re-entering a value through .Value which Excel would
# otherwise auto-coerce into a number.
for ($i = 0; $i -lt $sortedKeys.Count; $i++) {
    $srcRow = $sortedKeys[$i][1]
    $dstRow = $firstDataRow + $i
    $src = $ws.Range("A" + $srcRow + ":D" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":D" + $dstRow)
    $src.Cut($dst)
}
